# Add the author's affiliation paragraph right after the "Edison Achalma"
# author line (w:pStyle "Author") on the title page.
#
# Target result (per the diff): a new <w:p> with pStyle "Author" containing
# a single run "Escuela Profesional de Economía, Universidad Nacional de
# San Cristóbal de Huamanga", inserted immediately after the paragraph that
# just contains "Edison Achalma" in the Author style.

$d = $word.ActiveDocument

$affiliationText = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"

# Locate the specific "Edison Achalma" paragraph that uses the "Author"
# style (there are other "Edison Achalma" mentions further down the
# document, e.g. in "Nota de Autores", so match on style + exact text).
$target = $null
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text.Trim()
    if ($styleName -eq "Author" -and $text -eq "Edison Achalma") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Insert a new paragraph mark right after this paragraph's text.
    $insertionPoint = $target.Range.Duplicate
    $insertionPoint.Collapse(0)
    $insertionPoint.InsertAfter([char]13)

    # The freshly created paragraph is the next one; give it the "Author"
    # style and fill in the affiliation text.
    $newIndex = $target.Index + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Style = "Author"
    $newPara.Range.Text = $affiliationText
}
